$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98

# Row 4 updates
$ws.Range("G4").Value = 2.63
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 3.4

# Row 6 updates
$ws.Range("H6").Value = 3.55
$ws.Range("I6").Value = 7.9
$ws.Range("L6").Value = 7.3
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 6.6
$ws.Range("P6").Value = 2.72
$ws.Range("T6").Value = 2.5
$ws.Range("X6").Value = 6
$ws.Range("AC6").Value = 7.8
$ws.Range("AD6").Value = 7.2
$ws.Range("AE6").Value = 21
$ws.Range("AL6").Value = 120
$ws.Range("AT6").Value = 2.47
$ws.Range("AU6").Value = 8
$ws.Range("AV6").Value = 80
$ws.Range("AW6").Value = 8.75
